$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 99
$ws.Range("H99").Value = 570.7692
$ws.Range("I99").Value = 320
$ws.Range("J99").Value = 972
$ws.Range("K99").Value = 960
$ws.Range("L99").Value = 2916
$ws.Range("M99").Value = 538
$ws.Range("N99").Value = -5912
# Row 100
$ws.Range("H100").Value = 1698.3158
$ws.Range("I100").Value = 1443.6666
$ws.Range("J100").Value = 1815.8462
$ws.Range("K100").Value = 1443.6666
$ws.Range("L100").Value = 1815.8462
$ws.Range("M100").Value = -902.6666
$ws.Range("N100").Value = -2897.8462
# Row 101
$ws.Range("H101").Value = 8317.315000000001
$ws.Range("I101").Value = 530.5
$ws.Range("J101").Value = 13980.454
$ws.Range("K101").Value = 1591.5
$ws.Range("L101").Value = 41941.362
$ws.Range("M101").Value = 30.5
$ws.Range("N101").Value = -45185.362
# Row 103
$ws.Range("H103").Value = 945.9231
$ws.Range("J103").Value = 1199
$ws.Range("L103").Value = 3597
$ws.Range("N103").Value = -4769
# Row 106
$ws.Range("H106").Value = 5552
$ws.Range("I106").Value = 3268.3333
$ws.Range("K106").Value = 3268.3333
$ws.Range("M106").Value = -2637.3333
# Row 107
$ws.Range("H107").Value = 1026.3077
$ws.Range("I107").Value = 1095.4
$ws.Range("J107").Value = 796
$ws.Range("K107").Value = 1095.4
$ws.Range("L107").Value = 796
$ws.Range("M107").Value = 824.5999999999999
$ws.Range("N107").Value = -4636
# Row 137
$ws.Range("H137").Value = 30133.135
$ws.Range("I137").Value = 1356.84
$ws.Range("J137").Value = 90083.75
$ws.Range("K137").Value = 4070.52
$ws.Range("L137").Value = 270251.25
$ws.Range("M137").Value = -1520.52
$ws.Range("N137").Value = -275351.25

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1513.7142
$ws.Range("I61").Value = 1279.2307
$ws.Range("J61").Value = 2428.2
$ws.Range("K61").Value = 1279.2307
$ws.Range("L61").Value = 2428.2
$ws.Range("M61").Value = -1067.2307
$ws.Range("N61").Value = -2852.2
# Row 97
$ws.Range("H97").Value = 945.8570999999999
$ws.Range("I97").Value = 627.5
$ws.Range("J97").Value = 1370.3334
$ws.Range("K97").Value = 627.5
$ws.Range("L97").Value = 1370.3334
$ws.Range("M97").Value = -131.5
$ws.Range("N97").Value = -2362.3334
# Row 132
$ws.Range("H132").Value = 3576.4614
$ws.Range("I132").Value = 3923
$ws.Range("J132").Value = 3279.4285
$ws.Range("K132").Value = 11769
$ws.Range("L132").Value = 9838.2855
$ws.Range("M132").Value = -9239
$ws.Range("N132").Value = -14898.2855
# Row 136
$ws.Range("H136").Value = 1513.7142
$ws.Range("I136").Value = 1279.2307
$ws.Range("J136").Value = 2428.2
$ws.Range("K136").Value = 3837.6921
$ws.Range("L136").Value = 7284.599999999999
$ws.Range("M136").Value = -1287.6921
$ws.Range("N136").Value = -12384.6

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 350.88235
$ws.Range("I94").Value = 261.92307
$ws.Range("J94").Value = 640
$ws.Range("K94").Value = 261.92307
$ws.Range("L94").Value = 640
$ws.Range("M94").Value = 189.07693
$ws.Range("N94").Value = -1542
# Row 134
$ws.Range("H134").Value = 2045.6774
$ws.Range("I134").Value = 1843.6
$ws.Range("J134").Value = 2887.6667
$ws.Range("K134").Value = 5530.799999999999
$ws.Range("L134").Value = 8663.000100000001
$ws.Range("M134").Value = -2995.799999999999
$ws.Range("N134").Value = -13733.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 110
$ws.Range("H110").Value = 6679.8
$ws.Range("I110").Value = 6599.5
$ws.Range("J110").Value = 6733.3335
$ws.Range("K110").Value = 19798.5
$ws.Range("L110").Value = 20200.0005
$ws.Range("M110").Value = -15708.5
$ws.Range("N110").Value = -28380.0005
# Row 111
$ws.Range("H111").Value = 4000
$ws.Range("I111").Value = 5000
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 15000
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -11933
$ws.Range("N111").Value = -15134

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1064.9286
$ws.Range("I97").Value = 652.8570999999999
$ws.Range("K97").Value = 652.8570999999999
$ws.Range("M97").Value = -156.8570999999999
# Row 132
$ws.Range("H132").Value = 2621.7273
$ws.Range("I132").Value = 2646.973
$ws.Range("J132").Value = 2569.8333
$ws.Range("K132").Value = 7940.919
$ws.Range("L132").Value = 7709.499899999999
$ws.Range("M132").Value = -5410.919
$ws.Range("N132").Value = -12769.4999

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2070.258
$ws.Range("I7").Value = 1888.0741
$ws.Range("J7").Value = 3300
$ws.Range("K7").Value = 1888.0741
$ws.Range("L7").Value = 3300
$ws.Range("M7").Value = -1776.0741
$ws.Range("N7").Value = -3524
# Row 55
$ws.Range("H55").Value = 157
$ws.Range("I55").Value = 157
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 157
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()
# Row 100
$ws.Range("H100").Value = 85001230
$ws.Range("I100").Value = 3334130
$ws.Range("J100").Value = 166668340
$ws.Range("K100").Value = 3334130
$ws.Range("L100").Value = 166668340
$ws.Range("M100").Value = -3333589
$ws.Range("N100").Value = -166669422
# Row 126
$ws.Range("H126").Value = 2070.258
$ws.Range("I126").Value = 1888.0741
$ws.Range("J126").Value = 3300
$ws.Range("K126").Value = 5664.2223
$ws.Range("L126").Value = 9900
$ws.Range("M126").Value = -3194.2223
$ws.Range("N126").Value = -14840

$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 4000
$ws.Range("I14").Value = 5000
$ws.Range("J14").Value = 3000
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = -4832
$ws.Range("N14").Value = -3336
# Row 24
$ws.Range("H24").Value = 4933.3335
$ws.Range("I24").Value = 4900
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 4900
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = -4670
$ws.Range("N24").Value = -5460
# Row 40
$ws.Range("H40").Value = 20500
$ws.Range("J40").Value = 20500
$ws.Range("L40").Value = 20500
$ws.Range("N40").Value = -20798
# Row 135
$ws.Range("H135").Value = 38419.168
$ws.Range("J135").Value = 38419.168
$ws.Range("L135").Value = 38419.168
$ws.Range("N135").Value = -48559.168
